# feat(pipelining): add MOVS testing program #18
#
# 1) Rename the existing "Sheet1" to "Code".
# 2) Insert a new worksheet "Pipelining" right after "Code" and make it the
#    active sheet/tab.
# 3) Populate "Pipelining" with a small pipeline diagram (Fetch / Execute1 /
#    Execute2 / Store stages) for the four instructions already listed on
#    the "Code" sheet (rows 8-11, columns A:B).
# 4) Restore the selection on "Code" to the A8:B11 block (the source data
#    for the new diagram) and leave D4 selected on "Pipelining".

$wb = $excel.ActiveWorkbook

# --- 1) Rename Sheet1 -> Code -------------------------------------------
$code = $wb.Worksheets.Item(1)
$code.Name = "Code"

# --- 2) Add the new Pipelining sheet right after Code -------------------
$pipelining = $wb.Worksheets.Add([System.Type]::Missing, $code)
$pipelining.Name = "Pipelining"

# --- 3) Fill in the pipeline diagram ------------------------------------
$pipelining.Columns.Item(2).ColumnWidth = 21

# Row 1 - LOAD R1, [R0]
$pipelining.Range("A1").Value = "0x000"
$pipelining.Range("B1").Value = "LOAD R1, [R0]"
$pipelining.Range("C1").Value = "F"
$pipelining.Range("D1").Value = "E1"
$pipelining.Range("E1").Value = "E2"

# Row 2 - LOAD R2, [R0, #1]
$pipelining.Range("A2").Value = "0x001"
$pipelining.Range("B2").Value = "LOAD R2, [R0, #1]"
$pipelining.Range("D2").Value = "ST"
$pipelining.Range("E2").Value = "F"
$pipelining.Range("F2").Value = "E1"
$pipelining.Range("G2").Value = "E2"

# Row 3 - LOAD R3, [R0, #2]
$pipelining.Range("A3").Value = "0x002"
$pipelining.Range("B3").Value = "LOAD R3, [R0, #2]"
$pipelining.Range("F3").Value = "ST"
$pipelining.Range("G3").Value = "F"
$pipelining.Range("H3").Value = "E1"
$pipelining.Range("I3").Value = "E2"

# Row 4 - STP
$pipelining.Range("A4").Value = "0x003"
$pipelining.Range("B4").Value = "STP"
$pipelining.Range("H4").Value = "ST"
$pipelining.Range("I4").Value = "F"
$pipelining.Range("J4").Value = "E1"
$pipelining.Range("K4").Value = "E2"

# Left-align the instruction mnemonic column, matching the style already
# used for the same strings on the Code sheet.
$pipelining.Range("B1:B4").HorizontalAlignment = -4131

# --- 4) Selections -------------------------------------------------------
$code.Activate()
$code.Range("A8:B11").Select()

$pipelining.Activate()
$pipelining.Range("D4").Select()
